$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Helper: write a literal TEXT value into a cell without letting
# Excel's automatic number/percent parsing re-style the cell (typing
# "36.1%" directly would convert the cell to a Percentage number and
# mint/attach a new style index). We stage the text in a scratch cell
# (Z1, unused by this sheet) formatted as Text, copy it, then
# paste-special VALUES ONLY into the destination so the destination
# keeps its own existing style untouched.
# ------------------------------------------------------------------
function Set-LiteralText($rangeAddr, $text) {
    $scratch = $ws.Range("Z1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

# ------------------------------------------------------------------
# Helper: re-style a "session" row (columns A:I) by copying the
# formatting (fill/font/alignment) from a known-good row that already
# carries the desired style (the plain / non-pink "Recorded" look),
# then overwrite the Recorded-By / Students / Status cells with their
# new values. PasteSpecial(Formats) reuses the existing style index
# from the stylesheet instead of minting a new one.
# ------------------------------------------------------------------
function Set-RowRecorded($row, $recordedBy, $students) {
    $ws.Range("A4:I4").Copy()
    $ws.Range("A" + $row + ":I" + $row).PasteSpecial(-4122)
    $ws.Range("G" + $row).Value = $recordedBy
    $ws.Range("H" + $row).Value = $students
    $ws.Range("I" + $row).Value = "Recorded"
}

# ------------------------------------------------------------------
# Class Statistics block (K/L columns)
# ------------------------------------------------------------------
$ws.Range("L6").Value = 117          # Recorded Sessions
$ws.Range("L7").Value = 3            # Missing Sessions
Set-LiteralText "L9" "36.1%"         # Coverage %
Set-LiteralText "L10" "70.4%"        # Average Attendance %

# ------------------------------------------------------------------
# "Recorded By" cells: System now appears first, e.g.
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# ------------------------------------------------------------------
$swappedRecordedByRows = @(8, 9, 10, 35, 36, 37, 62, 63, 64, 89, 90, 91, `
    116, 117, 118, 143, 144, 145, 170, 197, 224, 251, 278, 305)
foreach ($r in $swappedRecordedByRows) {
    $ws.Range("G" + $r).Value = "System, dnasr281@gmail.com"
}

# ------------------------------------------------------------------
# Newly-recorded "session 11 / 14-12-2025" rows per group: these were
# "Not Recorded" (pink style) and are now "Recorded" (plain style),
# with a Recorded-By email and an actual attendance count.
# ------------------------------------------------------------------
Set-RowRecorded 12  "dnasr281@gmail.com" "19/26"   # B1A1
Set-RowRecorded 39  "dnasr281@gmail.com" "26/27"   # B1A2
Set-RowRecorded 66  "dnasr281@gmail.com" "17/26"   # B1B1
Set-RowRecorded 93  "dnasr281@gmail.com" "23/27"   # B1B2
Set-RowRecorded 120 "dnasr281@gmail.com" "24/30"   # B1C1
Set-RowRecorded 147 "dnasr281@gmail.com" "17/23"   # B1C2

# ------------------------------------------------------------------
# Per-group summary table (rows 15-20, columns O/P Recorded/Missing
# counts and R/S Coverage%/Avg Attendance% text).
# ------------------------------------------------------------------
$ws.Range("O15").Value = 10
$ws.Range("P15").Value = 1
Set-LiteralText "R15" "37.0%"
Set-LiteralText "S15" "80.8%"

$ws.Range("O16").Value = 11
$ws.Range("P16").Value = 0
Set-LiteralText "R16" "40.7%"
Set-LiteralText "S16" "77.4%"

$ws.Range("O17").Value = 11
$ws.Range("P17").Value = 0
Set-LiteralText "R17" "40.7%"
Set-LiteralText "S17" "63.6%"

$ws.Range("O18").Value = 11
$ws.Range("P18").Value = 0
Set-LiteralText "R18" "40.7%"
Set-LiteralText "S18" "67.0%"

$ws.Range("O19").Value = 11
$ws.Range("P19").Value = 0
Set-LiteralText "R19" "40.7%"
Set-LiteralText "S19" "69.4%"

$ws.Range("O20").Value = 10
$ws.Range("P20").Value = 1
Set-LiteralText "R20" "37.0%"
Set-LiteralText "S20" "73.0%"
